$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted above the existing data
# block (current row 90), pushing the existing rows 90-97 down to 91-98.
$ws.Rows("90:90").Insert()

# Populate the newly inserted row 90 with the new weekly record.
$ws.Range("A90").Value = 4
$ws.Range("B90").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C90").Value = "Los Lagos"
$ws.Range("D90").Value = 44449
$ws.Range("E90").Value = 10
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100102
$ws.Range("H90").Value = "Cítricos"
$ws.Range("I90").Value = 100102004
$ws.Range("J90").Value = "Mandarina"
$ws.Range("K90").Value = "Clementina"
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 600
$ws.Range("N90").Value = 6500
$ws.Range("O90").Value = 6500
$ws.Range("P90").Value = 6500
$ws.Range("Q90").Value = "$/bandeja 10 kilos"
$ws.Range("R90").Value = "Provincia de Limarí"
$ws.Range("S90").Value = 650
$ws.Range("T90").Value = 10
